$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sensor metadata rows to append (row 4, 5, 6)
# Columns: A=display_name, B=type, C=location_identifier, D=units, E=column_name, F=description

# Row 4 - cuteness sensor (entry order: E, B, C, D, A)
$ws.Range("E4").Value = "my_cuteness_sensor"
$ws.Range("B4").Value = "cute"
$ws.Range("C4").Value = "the heart"
$ws.Range("D4").Value = "kisses"
$ws.Range("A4").Value = "my cuteness sensor"

# Row 5 - coolness sensor (no location_identifier; entry order: B, D, E, A)
$ws.Range("B5").Value = "cool"
$ws.Range("D5").Value = "cigarettes"
$ws.Range("E5").Value = "my_coolness_sensor"
$ws.Range("A5").Value = "my coolness sensor"

# Row 6 - intelligence sensor (entry order: E, A, B, C, D)
$ws.Range("E6").Value = "intelligence_sensor"
$ws.Range("A6").Value = "my intelligence"
$ws.Range("B6").Value = "intl"
$ws.Range("C6").Value = "brain"
$ws.Range("D6").Value = "opions"

# Update the selection to reflect the new active cell D6
$ws.Range("D6").Select()
